# Append a new attendance record as row 5 (A1:E4 -> A1:E5), matching the
# existing rows' plain-text / numeric layout with no extra cell styling.
#
# Text cells are written via a quoted-string formula ( ="literal" ) instead
# of a plain .Value assignment: assigning the literal strings directly would
# let Excel's type-inference turn "10/05/2025" into a real date serial (and
# attach a date number-format style), which the source data does not use.
# The formula trick keeps the cells as plain text/"General" with no style,
# exactly like the pre-existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Formula = "=""10/05/2025"""
$ws.Range("B5").Value = 2702258535
$ws.Range("C5").Formula = "=""14:06:05"""
$ws.Range("D5").Formula = "="""""
$ws.Range("E5").Formula = "="""""
